$wb = $excel.ActiveWorkbook

$wsSite  = $wb.Worksheets.Item("Sitedeclaration")
$wsPages = $wb.Worksheets.Item("Pagesdeclarations")

# --- Pagesdeclarations: rewrite the PageName/URL list (dedupe fix + many new page paths) ---
$wsPages.Range("B2").Value  = "/news"
$wsPages.Range("B3").Value  = "/contact"
$wsPages.Range("B4").Value  = "/services/account-based-marketing-abm"
$wsPages.Range("B5").Value  = "/work"
$wsPages.Range("B6").Value  = "/about"
$wsPages.Range("B7").Value  = "/careers"
$wsPages.Range("B8").Value  = "/news"
$wsPages.Range("B9").Value  = "/contact"
$wsPages.Range("B10").Value = "/services/research"
$wsPages.Range("B11").Value = "/services/public-relations"
$wsPages.Range("B12").Value = "/services/creative-content"
$wsPages.Range("B13").Value = "/services/digital/user-experience-optimization"
$wsPages.Range("B14").Value = "/services/digital/digital-design"
$wsPages.Range("B15").Value = "/services/digital/website-development"
$wsPages.Range("B16").Value = "/services/digital/crm"
$wsPages.Range("B17").Value = "/services/account-based-marketing-abm/abm-strategy-and-programs"
$wsPages.Range("B18").Value = "/services/go-to-market/performance-marketing"
$wsPages.Range("B19").Value = "/about/leadership"
$wsPages.Range("B20").Value = "/locations"
$wsPages.Range("B21").Value = "/services/media/paid-search"

# B13 picked up formatting from a pasted source (dark grey "Aptos Narrow" text)
$wsPages.Range("B13").Font.Name = "Aptos Narrow"
$wsPages.Range("B13").Font.Color = 2368548

# --- Sitedeclaration: add a new (still blank) hyperlink-styled row under the table ---
$wsSite.Range("C3").Style = "Hyperlink"

$site_table = $wsSite.ListObjects.Item(1)
$site_table.Resize($wsSite.Range("A1:D3"))

# --- selection / active-sheet bookkeeping ---
$wsPages.Range("B20").Select()

$wsSite.Activate()
$wsSite.Range("D3").Select()
